# Trade #6 closed at 2026-02-16 22:52:32 - base_strategy DOWN +0.000%
# Append the new open trade as row 7 on both the "All Trades" and
# "base_strategy" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 7

    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = "2026-02-16"
    $ws.Cells.Item($row, 3).Value = "22:52:32"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 49.999998
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = 0
}
